$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column X. This shifts the existing
# "Snow.prec.2001.2030.mm.year" column (X) to Y, and the existing
# "Snow.prec.2071.2100.mm.year" column (Y) to Z.
$ws.Columns("X:X").Insert()

# New header for the inserted column
$ws.Range("X1").Value = "Current.snow.prec.offset.mm.year"

# New data for the inserted column (rows 2-23)
$values = @{
    2  = 139.699969852511
    3  = 130.451506091414
    4  = -187.879721584071
    5  = -624.370751352153
    6  = -13.2429241023732
    7  = -583.693898078582
    8  = 352.820918653345
    9  = 987.663930718751
    10 = 1041.80285004576
    11 = -711.694844660469
    12 = 265.490684529234
    13 = 92.5275662751825
    14 = 40.0608308308154
    15 = -24.8281128512736
    16 = 1083.27791288378
    17 = 4782.97014876069
    18 = 5312.51000096009
    19 = 177.616935774826
    20 = 547.994951903219
    21 = 361.206429474445
    22 = -228.849744831442
    23 = 185.918706696172
}

foreach ($row in $values.Keys) {
    $ws.Range("X$row").Value = $values[$row]
}
